$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price update: insert two rows for the newest reporting week at the top of
# this producer/variety block (row 64), pushing every later week down by two rows so
# the previously-last two weeks now land on the newly extended rows 112-113.
$ws.Range("A64:R65").Insert()

# Row 64: new week data
$ws.Range("A64").Value = 1
$ws.Range("B64").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C64").Value = "Arica y Parinacota"
$ws.Range("D64").Value = 44729
$ws.Range("E64").Value = 15
$ws.Range("F64").Value = 100112042
$ws.Range("G64").Value = "Locoto"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 140
$ws.Range("K64").Value = 21000
$ws.Range("L64").Value = 22000
$ws.Range("M64").Value = 21500
$ws.Range("N64").Value = "$/caja 20 kilos"
$ws.Range("O64").Value = "Región de Arica y Parinacota"
$ws.Range("P64").Value = 1075
$ws.Range("Q64").Value = 20
$ws.Range("R64").Value = "Hortaliza"

# Row 65: new week data
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44729
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100112042
$ws.Range("G65").Value = "Locoto"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Segunda"
$ws.Range("J65").Value = 140
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = 19500
$ws.Range("N65").Value = "$/caja 20 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 975
$ws.Range("Q65").Value = 20
$ws.Range("R65").Value = "Hortaliza"

